# Update "想去人数" (F column) figures for two events on both the
# "展览" sheet and the "全部类型" sheet.
#
#   展览      : F4 4830 -> 4842   (row for 南宁·2024良牙动漫秋季盛典)
#               F5   12 -> 15     (row for 南宁·花海演绎二次元水上派对)
#   全部类型  : F4 4830 -> 4842   (row for 南宁·2024良牙动漫秋季盛典)
#               F6   12 -> 15     (row for 南宁·花海演绎二次元水上派对)

$wb = $excel.ActiveWorkbook

$wsZhanLan = $wb.Worksheets.Item("展览")
$wsZhanLan.Range("F4").Value = 4842
$wsZhanLan.Range("F5").Value = 15

$wsQuanBu = $wb.Worksheets.Item("全部类型")
$wsQuanBu.Range("F4").Value = 4842
$wsQuanBu.Range("F6").Value = 15
